# MtomToCrss_Annual.xlsx - "May 2019 Official Run (Updated)" edit
#
# Every Trace* sheet's "ICS Credits.SystemEfficiencyICS_CA" column (E) is
# updated for the two modeled years (rows 2 and 3) to reflect the updated
# Mead Bank Parameters Table for the lower basin in the ICS Credits data
# object. A subset of traces also have their downstream PowellOperation
# totals (columns H/J, row 3) shift as a consequence of the updated Flaming
# Gorge rule / CRIT water-use schedule feeding into the same annual run.

$wb = $excel.ActiveWorkbook

$newE2 = 315218.56000000006
$newE3 = 377078.50320000004

# Row 3 columns H (PowellOperation.PowellWYRelease carry-through total) and
# J (PowellOperation.Compact Point Volume) that also shift for specific
# traces as a result of the Flaming Gorge / CRIT updates.
$H3J3 = @{
    4  = @(8539813.1611712184, 8571557.891171217)
    6  = @(12592234.043675648, 12804929.053675646)
    12 = @(9348251.9274542443, 9477419.7074542455)
    20 = @(11937964.933185101, 11999636.619185098)
    21 = @(10806309.720770467, 10934828.452770464)
    22 = @(11061200.617180195, 11343620.921180194)
    24 = @(12156576.581170088, 12566394.805170087)
    31 = @(12006616.126452897, 12154648.754452901)
    32 = @(10639120.274446048, 10864023.606446046)
    33 = @(13901991.423403217, 14193298.625403218)
    34 = @(12630518.852216702, 12904368.734216701)
    35 = @(10283454.571929816, 10494058.917929813)
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # "Sheet1" (the last tab) carries no trace data - nothing to update.
    if ($ws.Name -eq "Sheet1") {
        continue
    }

    $ws.Range("E2").Value = $newE2
    $ws.Range("E3").Value = $newE3

    if ($H3J3.ContainsKey($i)) {
        $vals = $H3J3[$i]
        $ws.Range("H3").Value = $vals[0]
        $ws.Range("J3").Value = $vals[1]
    }
}
